# Hortaliza, Vega Modelo de Temuco - Ciboulette: weekly update.
# Inserts one new data row (125) with a new price observation, shifting the
# existing rows 125-234 down to 126-235 (Excel's normal Rows.Insert behaviour).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 125, pushing rows 125..234 down to 126..235.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new weekly observation.
$ws.Cells.Item(125, 1).Value = 10
$ws.Cells.Item(125, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(125, 3).Value = "La Araucanía"
$ws.Cells.Item(125, 4).Value = 44658
$ws.Cells.Item(125, 5).Value = 9
$ws.Cells.Item(125, 6).Value = 100112039
$ws.Cells.Item(125, 7).Value = "Ciboulette"
$ws.Cells.Item(125, 8).Value = "Sin especificar"
$ws.Cells.Item(125, 9).Value = "Primera"
$ws.Cells.Item(125, 10).Value = 65
$ws.Cells.Item(125, 11).Value = 6000
$ws.Cells.Item(125, 12).Value = 6000
$ws.Cells.Item(125, 13).Value = 6000
$ws.Cells.Item(125, 14).Value = "`$/docena de atados"
$ws.Cells.Item(125, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(125, 16).Value = 2000
$ws.Cells.Item(125, 17).Value = 3
$ws.Cells.Item(125, 18).Value = "Hortaliza"
